$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'57.067.01"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "'2.981.10"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'500.04"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -5.18%  "
$ws.Range("D6").Value = "'136.86"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("D9").Value = "'7.27"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -5.41%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -4.85%  "
$ws.Range("D11").Value = "'0.356"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").Value = "'3.493.65"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").Value = "'0.0000159"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -6.36%  "
$ws.Range("D16").Value = "'57.139.30"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "'6.09"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "'2.984.36"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").Value = "'7.85"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").Value = "'319.81"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -5.66%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'5.72"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'0.492"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "'63.02"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "'0.0₃0889"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -9.51%  "
$ws.Range("D29").Value = "'6.61"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("D30").Value = "'7.09"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("D33").Value = "'20.13"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -5.15%  "
$ws.Range("D34").Value = "'154.73"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "'4.57"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "'5.77"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("E37").Value = "  -6.89%  "
$ws.Range("D38").Value = "'24.45"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("E39").Value = "  -5.94%  "
$ws.Range("D40").Value = "'37.80"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'3.015.04"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'3.73"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "'2.183.50"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -6.45%  "
$ws.Range("E46").Value = "  -6.79%  "
$ws.Range("D47").Value = "'5.95"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "'0.924"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -10.39%  "
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "'19.22"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -4.81%  "
$ws.Range("D51").Value = "'1.77"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -11.60%  "
